$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('F2').Value = 82
$ws.Range('H2').Value = 100
$ws.Range('I2').Value = 112
$ws.Range('E3').Value = 137
$ws.Range('H3').Value = 142
$ws.Range('E4').Value = 11
$ws.Range('D6').Value = 387
$ws.Range('E6').Value = 436
$ws.Range('F6').Value = 485
$ws.Range('G6').Value = 420
$ws.Range('I6').Value = 475
$ws.Range('D7').Value = 607
$ws.Range('E7').Value = 653
$ws.Range('F7').Value = 700
$ws.Range('G7').Value = 640
$ws.Range('H7').Value = 681
$ws.Range('I7').Value = 794

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('E6').Value = 48
$ws.Range('E7').Value = 61

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('E3').Value = 8
$ws.Range('E7').Value = 34

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('E4').Value = 4
$ws.Range('E5').Value = 8

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('E2').Value = 4
$ws.Range('E6').Value = 2
$ws.Range('E8').Value = 47
$ws.Range('F8').Value = 43
$ws.Range('G8').Value = 32
$ws.Range('F19').Value = 22
$ws.Range('G29').Value = 10
$ws.Range('E32').Value = 61
$ws.Range('E36').Value = 34
$ws.Range('E47').Value = 16
$ws.Range('E48').Value = 6
$ws.Range('D53').Value = 68
$ws.Range('E53').Value = 80
$ws.Range('H53').Value = 91
$ws.Range('I53').Value = 121
$ws.Range('I65').Value = 22
$ws.Range('D81').Value = 3
$ws.Range('E88').Value = 8
$ws.Range('D98').Value = 607
$ws.Range('E98').Value = 653
$ws.Range('F98').Value = 700
$ws.Range('G98').Value = 640
$ws.Range('H98').Value = 681
$ws.Range('I98').Value = 794

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('H2').Value = 12
$ws.Range('I2').Value = 12
$ws.Range('H3').Value = 17
$ws.Range('D6').Value = 40
$ws.Range('E6').Value = 62
$ws.Range('I6').Value = 77
$ws.Range('D7').Value = 68
$ws.Range('E7').Value = 80
$ws.Range('H7').Value = 91
$ws.Range('I7').Value = 121

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('D5').Value = 2
$ws.Range('D6').Value = 3

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I5').Value = 16
$ws.Range('I6').Value = 22

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('G5').Value = 8
$ws.Range('G6').Value = 10

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('F5').Value = 15
$ws.Range('F6').Value = 22

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('E4').Value = 2
$ws.Range('E6').Value = 4

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('E5').Value = 5
$ws.Range('E6').Value = 6

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('E5').Value = 10
$ws.Range('E6').Value = 16

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('F2').Value = 8
$ws.Range('E5').Value = 37
$ws.Range('G5').Value = 23
$ws.Range('E6').Value = 47
$ws.Range('F6').Value = 43
$ws.Range('G6').Value = 32

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('E4').Value = 2
$ws.Range('E5').Value = 2
